$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 1
    5  = -4
    6  = -1
    8  = -2
    9  = -1
    10 = -5
    11 = -4
    12 = 1
    13 = -5
    14 = -1
    15 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
